$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 995
$ws.Range("B2").Value = 995
$ws.Range("C2").Value = 995
$ws.Range("D2").Value = 995
$ws.Range("G2").Value = 995
